# fix(gui) step 1 and 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Bump the list's date by one day
$ws.Range("A1").Value = 45309

# Step 1 prices ("COMUN" table, rows 35-39)
$ws.Range("D35").Value = 3030.119
$ws.Range("D36").Value = 3310
$ws.Range("D37").Value = 3310
$ws.Range("D38").Value = 3310
$ws.Range("D39").Value = 8628.17

# Step 2 prices ("CON TOPE" table, rows 42-46)
$ws.Range("D42").Value = 3267.016
$ws.Range("D43").Value = 3500
$ws.Range("D44").Value = 3500
$ws.Range("D45").Value = 3500
$ws.Range("D46").Value = 9649.955
